$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# Row 2 - Bitcoin
Set-TextValue "D2" "30.768.68"
Set-TextValue "E2" "  +0.69%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.949.38"
Set-TextValue "E3" "  +1.59%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "247.54"
Set-TextValue "E5" "  +1.12%  "

# Row 6 - USDC
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.06%  "

# Row 7 - XRP
Set-TextValue "D7" "0.4804"
Set-TextValue "E7" "  -1.50%  "

# Row 8 - Cardano
Set-TextValue "E8" "  +1.55%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.06821"
Set-TextValue "E9" "  +1.28%  "

# Row 10 - Litecoin
Set-TextValue "D10" "113.12"
Set-TextValue "E10" "  +5.55%  "

# Row 11 - Solana
Set-TextValue "D11" "19.53"
Set-TextValue "E11" "  +4.17%  "

# Row 12 - WrappedEther
Set-TextValue "D12" "1.949.78"
Set-TextValue "E12" "  +1.61%  "

# Row 13 - Polkadot
Set-TextValue "D13" "5.575"
Set-TextValue "E13" "  +5.14%  "

# Row 14 - TRON
Set-TextValue "D14" "0.07648"
Set-TextValue "E14" "  -0.14%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.6909"
Set-TextValue "E15" "  +3.25%  "

# Row 16 - BitcoinCash
Set-TextValue "D16" "298.53"
Set-TextValue "E16" "  +7.00%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "30.748.66"
Set-TextValue "E17" "  +0.62%  "

# Row 18 - Avalanche
Set-TextValue "D18" "13.29"
Set-TextValue "E18" "  +3.65%  "

# Row 19 - was ShibaInu, now Uniswap
Set-TextValue "B19" "Uniswap"
Set-TextValue "C19" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D19" "5.664"
Set-TextValue "E19" "  +2.85%  "

# Row 20 - was Uniswap, now ShibaInu
Set-TextValue "B20" "ShibaInu"
Set-TextValue "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.000007708"
Set-TextValue "E20" "  +1.87%  "

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue "D21" "2.198.49"
Set-TextValue "E21" "  +1.39%  "

# Row 22 - Dai
Set-TextValue "D22" "0.9998"
Set-TextValue "E22" "  -0.06%  "

# Row 24 - Chainlink
Set-TextValue "D24" "6.584"
Set-TextValue "E24" "  +2.37%  "

# Row 25 - Cosmos
Set-TextValue "D25" "9.744"
Set-TextValue "E25" "  +3.03%  "

# Row 26 - Monero
Set-TextValue "D26" "167.66"
Set-TextValue "E26" "  +1.85%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "20.51"
Set-TextValue "E27" "  +1.53%  "

# Row 28 - LidoDAOToken
Set-TextValue "D28" "2.177"
Set-TextValue "E28" "  +3.45%  "

# Row 29 - Stellar
Set-TextValue "D29" "0.1088"
Set-TextValue "E29" "  +3.08%  "

# Row 30 - Toncoin
Set-TextValue "D30" "1.431"
Set-TextValue "E30" "  +1.79%  "

# Row 31 - Filecoin
Set-TextValue "D31" "4.526"
Set-TextValue "E31" "  +11.77%  "

# Row 32 - InternetComputer(DFINITY)
Set-TextValue "D32" "4.421"
Set-TextValue "E32" "  +6.34%  "

# Row 33 - Hedera
Set-TextValue "D33" "0.05060"
Set-TextValue "E33" "  +1.10%  "

# Row 34 - ImmutableX
Set-TextValue "D34" "0.7790"
Set-TextValue "E34" "  +6.31%  "

# Row 35 - ARBITRUM
Set-TextValue "E35" "  +2.31%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.02073"
Set-TextValue "E36" "  +1.66%  "

# Row 37 - HuobiToken
Set-TextValue "D37" "2.730"
Set-TextValue "E37" "  -0.05%  "

# Row 38 - MXToken
Set-TextValue "D38" "2.707"
Set-TextValue "E38" "  +1.07%  "

# Row 39 - was Quant, now RenderToken
Set-TextValue "B39" "RenderToken"
Set-TextValue "C39" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D39" "2.042"
Set-TextValue "E39" "  +1.17%  "

# Row 40 - was RenderToken, now Quant
Set-TextValue "B40" "Quant"
Set-TextValue "C40" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D40" "111.24"
Set-TextValue "E40" "  -0.38%  "

# Row 41 - TheSandbox
Set-TextValue "D41" "0.4474"
Set-TextValue "E41" "  +0.45%  "

# Row 42 - TrustWalletToken
Set-TextValue "D42" "0.8759"
Set-TextValue "E42" "  +0.46%  "

# Row 43 - FraxShare
Set-TextValue "D43" "5.967"
Set-TextValue "E43" "  +1.02%  "

# Row 44 - Aave
Set-TextValue "D44" "71.75"
Set-TextValue "E44" "  +5.48%  "

# Row 45 - PaxDollar
Set-TextValue "D45" "0.9996"
Set-TextValue "E45" "  -0.12%  "

# Row 46 - Aptos
Set-TextValue "D46" "7.408"
Set-TextValue "E46" "  +1.81%  "

# Row 47 - EnergySwap
Set-TextValue "D47" "9.513"
Set-TextValue "E47" "  +2.42%  "

# Row 48 - BitcoinSV
Set-TextValue "D48" "49.09"
Set-TextValue "E48" "  +0.66%  "

# Row 49 - Algorand
Set-TextValue "D49" "0.1257"
Set-TextValue "E49" "  +0.19%  "

# Row 50 - was Elrond, now WOONetwork
Set-TextValue "B50" "WOONetwork"
Set-TextValue "C50" "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue "D50" "0.2557"
Set-TextValue "E50" "  +2.31%  "

# Row 51 - was WOONetwork, now Elrond
Set-TextValue "B51" "Elrond"
Set-TextValue "C51" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D51" "35.53"
Set-TextValue "E51" "  +1.96%  "
